# score_swap_column: the "midterm1"/"midterm2" column headers were
# swapped by mistake (column B was labeled "midterm2" and column C was
# labeled "midterm1", while the underlying score data was actually for
# midterm1 in column B and midterm2 in column C). Correct this by
# swapping just the header labels in B1 and C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerB = $ws.Range("B1").Value()
$headerC = $ws.Range("C1").Value()

$ws.Range("B1").Value = $headerC
$ws.Range("C1").Value = $headerB
